# Apply the edit described by the commit:
#  - swap the longitude/latitude header cells on "Warehouse" (B1/C1) and
#    "Transportation" (D1/E1)
#  - drop the stray empty formatted cell J8 on "Transportation"
#  - update each sheet's remembered selection (activeCell/sqref)
#  - make "Warehouse" the active/selected tab again (it was "Region")

$wb = $excel.ActiveWorkbook

$wsWarehouse = $wb.Worksheets.Item("Warehouse")
$wsRegion = $wb.Worksheets.Item("Region")
$wsTransportation = $wb.Worksheets.Item("Transportation")

# --- Region: just move the remembered selection, and make sure it is no
#     longer flagged as the active tab.
$wsRegion.Activate()
$wsRegion.Range("O21").Select()

# --- Transportation: move selection, swap D1/E1 header labels, and clear
#     the stray J8 cell.
$wsTransportation.Activate()
$wsTransportation.Range("E6").Select()

$d1 = $wsTransportation.Range("D1").Value2
$e1 = $wsTransportation.Range("E1").Value2
$wsTransportation.Range("D1").Value = $e1
$wsTransportation.Range("E1").Value = $d1

$wsTransportation.Range("J8").Clear()

# --- Warehouse: move selection, swap B1/C1 header labels, and re-activate
#     it so it becomes the workbook's selected tab (matching the diff,
#     which drops bookViews/workbookView@activeTab back to the default).
$wsWarehouse.Activate()
$wsWarehouse.Range("G12").Select()

$b1 = $wsWarehouse.Range("B1").Value2
$c1 = $wsWarehouse.Range("C1").Value2
$wsWarehouse.Range("B1").Value = $c1
$wsWarehouse.Range("C1").Value = $b1
